$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New IP addresses for each data row (A3:A38), derived from the commit's
# "IP addresses updated" change.
$newIPs = @{
    3  = "129.62.150.35"
    4  = "129.62.150.36"
    5  = "129.62.150.37"
    6  = "129.62.150.38"
    7  = "129.62.150.39"
    8  = "129.62.150.40"
    9  = "129.62.150.41"
    10 = "129.62.150.42"
    11 = "129.62.150.43"
    12 = "129.62.150.44"
    13 = "129.62.150.45"
    14 = "129.62.150.46"
    15 = "129.62.150.23"
    16 = "129.62.150.24"
    17 = "129.62.150.25"
    18 = "129.62.150.26"
    19 = "129.62.150.27"
    20 = "129.62.150.28"
    21 = "129.62.150.29"
    22 = "129.62.150.30"
    23 = "129.62.150.31"
    24 = "129.62.150.32"
    25 = "129.62.150.33"
    26 = "129.62.150.34"
    27 = "129.62.150.11"
    28 = "129.62.150.12"
    29 = "129.62.150.13"
    30 = "129.62.150.14"
    31 = "129.62.150.15"
    32 = "129.62.150.16"
    33 = "129.62.150.17"
    34 = "129.62.150.18"
    35 = "129.62.150.19"
    36 = "129.62.150.20"
    37 = "129.62.150.21"
    38 = "129.62.150.22"
}

$writeOrder = @(27,28,29,30,31,32,33,34,35,36,37,38, `
                15,16,17,18,19,20,21,22,23,24,25,26, `
                3,4,5,6,7,8,9,10,11,12,13,14)
foreach ($row in $writeOrder) {
    $ws.Cells.Item($row, 1).Value = $newIPs[$row]
}

# Column A got a bit wider to fit the new (longer) IP values.
$ws.Columns.Item(1).ColumnWidth = 12

# Scroll/selection moved further down the sheet.
$excel.ActiveWindow.ScrollRow = 32
$ws.Range("E36").Select()
